# The deck originally carried two themes in its package:
#   ppt/theme/theme1.xml  -> "Integral"      (used by the real slide master -> what
#                                              the audience actually sees)
#   ppt/theme/theme2.xml  -> "Office Theme"  (used only by the notes master)
#
# The authored change swaps which theme backs the visible presentation: the
# slide master ends up on the stock "Office Theme" colour values while the
# "Integral" palette is pushed onto the (otherwise unreachable from the
# object model) notes-only theme slot.
#
# The PowerPoint object model doesn't expose a "swap the two theme parts"
# verb, but it does expose the live colour slots of the active theme via
# Slide.ThemeColorScheme (1=dk1 .. 12=folHlink), and writes there land in
# ppt/theme/theme1.xml - the part that actually drives the deck's look.
# dk1/lt1 are already identical between the two themes, so only the other
# ten slots need to move to their "Office Theme" values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (slot, target RGB as packed 0xBBGGRR used by PowerPoint's RGB())
$tcs.Item(3).RGB  = 6968388    # dk2       44546A
$tcs.Item(4).RGB  = 15132391   # lt2       E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1   5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2   ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3   A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4   FFC000
$tcs.Item(9).RGB  = 12874308   # accent5   4472C4
$tcs.Item(10).RGB = 4697456    # accent6   70AD47
$tcs.Item(11).RGB = 12673797   # hlink     0563C1
$tcs.Item(12).RGB = 7491477    # folHlink  954F72
